$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing two columns (J:K) -- the "Emerging" group shrinks from
# 5 columns (G:K) to 4 (F:I), matching the "Developed" group (B:E).
$ws.Range("J1:K4").EntireColumn.Delete()

# Re-split the two header groups: was B1:F1 / G1:K1, now B1:E1 / F1:I1.
$ws.Range("B1:F1").UnMerge()
$ws.Range("G1:I1").UnMerge()

# The "Emerging" label moves from G1 (old second-group anchor) to F1 (new one).
$ws.Range("F1").Value = $ws.Range("G1").Value()
$ws.Range("G1").Value = ""

$ws.Range("B1:E1").Merge()
$ws.Range("F1:I1").Merge()

# Merge() redraws a box border around the merged area; restore the sheet's
# uniform thin-border-all-sides header style so it matches the rest of row 1.
$ws.Range("A1:I1").Borders.LineStyle = 1

# Row 2 quartile labels: each group now only spans 0..3 instead of 0..4.
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 3

# Row 4: new forward-return figures from the re-run factor processing.
$ws.Range("B4").Value = 0.008933448839533037
$ws.Range("C4").Value = 0.008794247992597764
$ws.Range("D4").Value = 0.00702469826538205
$ws.Range("E4").Value = 0.006055623174030418
$ws.Range("F4").Value = 0.01626836458133035
$ws.Range("G4").Value = 0.007863423373416115
$ws.Range("H4").Value = 0.01058097490108159
$ws.Range("I4").Value = 0.01166751238240992

Write-Output "done"
